# Update PLC data 2025-10-13 13:55:00
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 166575
$ws.Range("C4").Value = 157475
$ws.Range("C5").Value = 9100
$ws.Range("C8").Value = 65.29000000000001
